$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 2).Value = 14.63447178262513
$ws.Cells.Item(2, 3).Value = 9.928194497479533
$ws.Cells.Item(2, 4).Value = 14.19146109943217
$ws.Cells.Item(2, 5).Value = 14.99917169517252
$ws.Cells.Item(2, 7).Value = 43.71838808127639
$ws.Cells.Item(2, 8).Value = 18.1677928269181
$ws.Cells.Item(2, 10).Value = 8.938778278569178
$ws.Cells.Item(2, 11).Value = 10.43373911065065
$ws.Cells.Item(2, 12).Value = 11.73334175441157
$ws.Cells.Item(2, 15).Value = 29.61596429947628
# Row 3
$ws.Cells.Item(3, 2).Value = 14.43175023206108
$ws.Cells.Item(3, 3).Value = 9.936574106402505
$ws.Cells.Item(3, 4).Value = 14.18724944645671
$ws.Cells.Item(3, 5).Value = 15.0240717915005
$ws.Cells.Item(3, 7).Value = 43.84782588405859
$ws.Cells.Item(3, 8).Value = 18.22164839699431
$ws.Cells.Item(3, 10).Value = 8.952127706943084
$ws.Cells.Item(3, 11).Value = 10.27966833475122
$ws.Cells.Item(3, 12).Value = 11.73319092111807
$ws.Cells.Item(3, 15).Value = 29.70881781553483
# Row 4
$ws.Cells.Item(4, 2).Value = 14.30824148186578
$ws.Cells.Item(4, 3).Value = 9.942264362760685
$ws.Cells.Item(4, 4).Value = 14.18720349660893
$ws.Cells.Item(4, 5).Value = 15.04142807395589
$ws.Cells.Item(4, 7).Value = 43.93731917901619
$ws.Cells.Item(4, 8).Value = 18.25721671696111
$ws.Cells.Item(4, 10).Value = 8.960783423272202
$ws.Cells.Item(4, 11).Value = 10.18529639603996
$ws.Cells.Item(4, 12).Value = 11.73458469076688
$ws.Cells.Item(4, 15).Value = 29.77093809211249
# Row 5
$ws.Cells.Item(5, 2).Value = 14.25821120749323
$ws.Cells.Item(5, 3).Value = 9.944720707849934
$ws.Cells.Item(5, 4).Value = 14.18782520034126
$ws.Cells.Item(5, 5).Value = 15.04902113806031
$ws.Cells.Item(5, 7).Value = 43.97630143544674
$ws.Cells.Item(5, 8).Value = 18.27234027014256
$ws.Cells.Item(5, 10).Value = 8.96442647763371
$ws.Cells.Item(5, 11).Value = 10.14694049313715
$ws.Cells.Item(5, 12).Value = 11.73552762429969
$ws.Cells.Item(5, 15).Value = 29.79753611658895
# Row 6
$ws.Cells.Item(6, 2).Value = 14.24992357204987
$ws.Cells.Item(6, 3).Value = 9.945136900810663
$ws.Cells.Item(6, 4).Value = 14.18796715834367
$ws.Cells.Item(6, 5).Value = 15.05031339074709
$ws.Cells.Item(6, 7).Value = 43.98292600460606
$ws.Cells.Item(6, 8).Value = 18.27488953503852
$ws.Cells.Item(6, 10).Value = 8.965038406613845
$ws.Cells.Item(6, 11).Value = 10.14057894204765
$ws.Cells.Item(6, 12).Value = 11.73570687734838
$ws.Cells.Item(6, 15).Value = 29.80203019210251
# Row 7
$ws.Cells.Item(7, 2).Value = 14.30756546224289
$ws.Cells.Item(7, 3).Value = 9.942296932494859
$ws.Cells.Item(7, 4).Value = 14.18720928600052
$ws.Cells.Item(7, 5).Value = 15.04152836982746
$ws.Cells.Item(7, 7).Value = 43.93783474058847
$ws.Cells.Item(7, 8).Value = 18.25741813055025
$ws.Cells.Item(7, 10).Value = 8.960832085547654
$ws.Cells.Item(7, 11).Value = 10.18477864588522
$ws.Cells.Item(7, 12).Value = 11.73459588778825
$ws.Cells.Item(7, 15).Value = 29.77129160672457
# Row 8
$ws.Cells.Item(8, 2).Value = 14.56440854842272
$ws.Cells.Item(8, 3).Value = 9.930970904919475
$ws.Cells.Item(8, 4).Value = 14.18948285890991
$ws.Cells.Item(8, 5).Value = 15.00732835308371
$ws.Cells.Item(8, 7).Value = 43.76093544991109
$ws.Cells.Item(8, 8).Value = 18.18584339068655
$ws.Cells.Item(8, 10).Value = 8.943286072499133
$ws.Cells.Item(8, 11).Value = 10.38059450644258
$ws.Cells.Item(8, 12).Value = 11.73298199640015
$ws.Cells.Item(8, 15).Value = 29.64691952846857
# Row 9
$ws.Cells.Item(9, 2).Value = 15.0731572828234
$ws.Cells.Item(9, 3).Value = 9.913066982101178
$ws.Cells.Item(9, 4).Value = 14.21400820758805
$ws.Cells.Item(9, 5).Value = 14.95665096959937
$ws.Cells.Item(9, 7).Value = 43.49378362509564
$ws.Cells.Item(9, 8).Value = 18.06531318509347
$ws.Cells.Item(9, 10).Value = 8.912506147825068
$ws.Cells.Item(9, 11).Value = 10.76448275102003
$ws.Cells.Item(9, 12).Value = 11.7415551661359
$ws.Cells.Item(9, 15).Value = 29.44359407584708
# Row 10
$ws.Cells.Item(10, 2).Value = 15.44663421344819
$ws.Cells.Item(10, 3).Value = 9.902512897541712
$ws.Cells.Item(10, 4).Value = 14.24412656166797
$ws.Cells.Item(10, 5).Value = 14.9293859843638
$ws.Cells.Item(10, 7).Value = 43.34645770939375
$ws.Cells.Item(10, 8).Value = 17.98882647447594
$ws.Cells.Item(10, 10).Value = 8.892082627783457
$ws.Cells.Item(10, 11).Value = 11.04394236622549
$ws.Cells.Item(10, 12).Value = 11.75493019400569
$ws.Cells.Item(10, 15).Value = 29.31899452921864
# Row 11
$ws.Cells.Item(11, 2).Value = 15.61577142107846
$ws.Cells.Item(11, 3).Value = 9.898270779055826
$ws.Cells.Item(11, 4).Value = 14.26042050593385
$ws.Cells.Item(11, 5).Value = 14.91914082007314
$ws.Cells.Item(11, 7).Value = 43.29012801500895
$ws.Cells.Item(11, 8).Value = 17.95664664176698
$ws.Cells.Item(11, 10).Value = 8.883262608918512
$ws.Cells.Item(11, 11).Value = 11.17000105117365
$ws.Cells.Item(11, 12).Value = 11.76253223131589
$ws.Cells.Item(11, 15).Value = 29.26770317307262
# Row 12
$ws.Cells.Item(12, 2).Value = 15.67965339926558
$ws.Cells.Item(12, 3).Value = 9.896744341424149
$ws.Cells.Item(12, 4).Value = 14.26695993783791
$ws.Cells.Item(12, 5).Value = 14.91557092134673
$ws.Cells.Item(12, 7).Value = 43.27033906709002
$ws.Cells.Item(12, 8).Value = 17.94483663716465
$ws.Cells.Item(12, 10).Value = 8.87999005074585
$ws.Cells.Item(12, 11).Value = 11.21754136646958
$ws.Cells.Item(12, 12).Value = 11.76562731961491
$ws.Cells.Item(12, 15).Value = 29.24905636420261
# Row 13
$ws.Cells.Item(13, 2).Value = 15.66590357672171
$ws.Cells.Item(13, 3).Value = 9.89706953750694
$ws.Cells.Item(13, 4).Value = 14.2655351915039
$ws.Cells.Item(13, 5).Value = 14.91632599764604
$ws.Cells.Item(13, 7).Value = 43.27453232515819
$ws.Cells.Item(13, 8).Value = 17.94736342138363
$ws.Cells.Item(13, 10).Value = 8.880691862204371
$ws.Cells.Item(13, 11).Value = 11.20731204448222
$ws.Cells.Item(13, 12).Value = 11.76495114536874
$ws.Cells.Item(13, 15).Value = 29.25303775328158
# Row 14
$ws.Cells.Item(14, 2).Value = 15.62103069128415
$ws.Cells.Item(14, 3).Value = 9.898143597992865
$ws.Cells.Item(14, 4).Value = 14.26095112996794
$ws.Cells.Item(14, 5).Value = 14.91884091853544
$ws.Cells.Item(14, 7).Value = 43.28846903930239
$ws.Cells.Item(14, 8).Value = 17.95566749479277
$ws.Cells.Item(14, 10).Value = 8.882992024673786
$ws.Cells.Item(14, 11).Value = 11.17391637002464
$ws.Cells.Item(14, 12).Value = 11.7627825410517
$ws.Cells.Item(14, 15).Value = 29.26615352828634
# Row 15
$ws.Cells.Item(15, 2).Value = 15.59352138302078
$ws.Cells.Item(15, 3).Value = 9.898811891799028
$ws.Cells.Item(15, 4).Value = 14.25819123743651
$ws.Cells.Item(15, 5).Value = 14.92042169727089
$ws.Cells.Item(15, 7).Value = 43.29720661024699
$ws.Cells.Item(15, 8).Value = 17.9608029156067
$ws.Cells.Item(15, 10).Value = 8.884409708093168
$ws.Cells.Item(15, 11).Value = 11.15343390259951
$ws.Cells.Item(15, 12).Value = 11.76148232725392
$ws.Cells.Item(15, 15).Value = 29.27428842392307
# Row 16
$ws.Cells.Item(16, 2).Value = 15.43556046944319
$ws.Cells.Item(16, 3).Value = 9.902801342193476
$ws.Cells.Item(16, 4).Value = 14.24311359103325
$ws.Cells.Item(16, 5).Value = 14.93009890354571
$ws.Cells.Item(16, 7).Value = 43.35035450511906
$ws.Cells.Item(16, 8).Value = 17.99098209755309
$ws.Cells.Item(16, 10).Value = 8.892668483542501
$ws.Cells.Item(16, 11).Value = 11.03567902782209
$ws.Cells.Item(16, 12).Value = 11.75446376460044
$ws.Cells.Item(16, 15).Value = 29.32245508734778
# Row 17
$ws.Cells.Item(17, 2).Value = 15.33842065949176
$ws.Cells.Item(17, 3).Value = 9.905391610905511
$ws.Cells.Item(17, 4).Value = 14.23452560602441
$ws.Cells.Item(17, 5).Value = 14.93658783983712
$ws.Cells.Item(17, 7).Value = 43.38570040065523
$ws.Cells.Item(17, 8).Value = 18.01016554538786
$ws.Cells.Item(17, 10).Value = 8.897855328507022
$ws.Cells.Item(17, 11).Value = 10.96313631333152
$ws.Cells.Item(17, 12).Value = 11.75054551715986
$ws.Cells.Item(17, 15).Value = 29.35338497674515
# Row 18
$ws.Cells.Item(18, 2).Value = 15.28248026715922
$ws.Cells.Item(18, 3).Value = 9.906934104619786
$ws.Cells.Item(18, 4).Value = 14.22983032103016
$ws.Cells.Item(18, 5).Value = 14.94052324605279
$ws.Cells.Item(18, 7).Value = 43.40703632772644
$ws.Cells.Item(18, 8).Value = 18.02144544087529
$ws.Cells.Item(18, 10).Value = 8.900882990894331
$ws.Cells.Item(18, 11).Value = 10.92131345039881
$ws.Cells.Item(18, 12).Value = 11.74843481497094
$ws.Cells.Item(18, 15).Value = 29.37168226326272
# Row 19
$ws.Cells.Item(19, 2).Value = 15.26352988716548
$ws.Cells.Item(19, 3).Value = 9.907465420250814
$ws.Cells.Item(19, 4).Value = 14.22828264065913
$ws.Cells.Item(19, 5).Value = 14.9418906130368
$ws.Cells.Item(19, 7).Value = 43.41443291878259
$ws.Cells.Item(19, 8).Value = 18.02530688917299
$ws.Cells.Item(19, 10).Value = 8.901915727656146
$ws.Cells.Item(19, 11).Value = 10.90713737496236
$ws.Cells.Item(19, 12).Value = 11.74774477799953
$ws.Cells.Item(19, 15).Value = 29.37796449355386
# Row 20
$ws.Cells.Item(20, 2).Value = 15.34876883239376
$ws.Cells.Item(20, 3).Value = 9.905110427426296
$ws.Cells.Item(20, 4).Value = 14.23541455290529
$ws.Cells.Item(20, 5).Value = 14.93587606100858
$ws.Cells.Item(20, 7).Value = 43.38183362672851
$ws.Cells.Item(20, 8).Value = 18.00809796690822
$ws.Cells.Item(20, 10).Value = 8.897298594581688
$ws.Cells.Item(20, 11).Value = 10.97086908059588
$ws.Cells.Item(20, 12).Value = 11.75094783764756
$ws.Cells.Item(20, 15).Value = 29.35003993247945
# Row 21
$ws.Cells.Item(21, 2).Value = 15.63421591562857
$ws.Cells.Item(21, 3).Value = 9.897825953637298
$ws.Cells.Item(21, 4).Value = 14.26228758835443
$ws.Cells.Item(21, 5).Value = 14.91809382437088
$ws.Cells.Item(21, 7).Value = 43.28433360690108
$ws.Cells.Item(21, 8).Value = 17.95321818915218
$ws.Cells.Item(21, 10).Value = 8.882314584595537
$ws.Cells.Item(21, 11).Value = 11.18373111557546
$ws.Cells.Item(21, 12).Value = 11.76341365625605
$ws.Cells.Item(21, 15).Value = 29.2622800352128
# Row 22
$ws.Cells.Item(22, 2).Value = 15.81977567507916
$ws.Cells.Item(22, 3).Value = 9.893531005636055
$ws.Cells.Item(22, 4).Value = 14.282001158399
$ws.Cells.Item(22, 5).Value = 14.9082770963448
$ws.Cells.Item(22, 7).Value = 43.22960084162707
$ws.Cells.Item(22, 8).Value = 17.91954146692876
$ws.Cells.Item(22, 10).Value = 8.872914339782486
$ws.Cells.Item(22, 11).Value = 11.32169181592381
$ws.Cells.Item(22, 12).Value = 11.77282095758682
$ws.Cells.Item(22, 15).Value = 29.2094481325414
# Row 23
$ws.Cells.Item(23, 2).Value = 15.72084858828348
$ws.Cells.Item(23, 3).Value = 9.895780811339232
$ws.Cells.Item(23, 4).Value = 14.27128415571118
$ws.Cells.Item(23, 5).Value = 14.91335151631986
$ws.Cells.Item(23, 7).Value = 43.25798880482164
$ws.Cells.Item(23, 8).Value = 17.93731499275199
$ws.Cells.Item(23, 10).Value = 8.877895596499028
$ws.Cells.Item(23, 11).Value = 11.24817879702915
$ws.Cells.Item(23, 12).Value = 11.76768544887658
$ws.Cells.Item(23, 15).Value = 29.2372311933181
# Row 24
$ws.Cells.Item(24, 2).Value = 15.34409071001807
$ws.Cells.Item(24, 3).Value = 9.905237384377795
$ws.Cells.Item(24, 4).Value = 14.23501190576414
$ws.Cells.Item(24, 5).Value = 14.93619721815399
$ws.Cells.Item(24, 7).Value = 43.38357863358932
$ws.Cells.Item(24, 8).Value = 18.00903193703022
$ws.Cells.Item(24, 10).Value = 8.89755015168938
$ws.Cells.Item(24, 11).Value = 10.96737345803683
$ws.Cells.Item(24, 12).Value = 11.75076550625783
$ws.Cells.Item(24, 15).Value = 29.35155062198901
# Row 25
$ws.Cells.Item(25, 2).Value = 14.93533740093401
$ws.Cells.Item(25, 3).Value = 9.91745206528401
$ws.Cells.Item(25, 4).Value = 14.20523767910848
$ws.Cells.Item(25, 5).Value = 14.96860820073593
$ws.Cells.Item(25, 7).Value = 43.55748164768766
$ws.Cells.Item(25, 8).Value = 18.09579926777751
$ws.Cells.Item(25, 10).Value = 8.92044674653804
$ws.Cells.Item(25, 11).Value = 10.66090577214694
$ws.Cells.Item(25, 12).Value = 11.73798700939521
$ws.Cells.Item(25, 15).Value = 29.49425037770871
